# CC_TA_OA_TEMPLATE.xlsx edit
# 1) NB Reporting: Capital calls exporting - replace hard-coded "Additional subscription"
#    use-of-funds option with a parametrised placeholder.
# 2) PE: M2S2 memo adjustments - Firm and Fund profile: replace the hard-coded
#    director name and the Russian visa / "doer" (executor) names with
#    parametrised placeholders so the export service can fill them in per document.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form")

# --- Use of Funds dropdown option (row 4) ---
$ws.Range("D4").Value = "<investment_category>"

# --- Director signature line (row 31) ---
$ws.Range("C31").Value = "<DIRECTORNAME>-Director"

# --- Visa / sign-off block (rows 82-84) and "doer" (executor) line (row 89) ---
$ws.Range("C82").Value = "<viza_1>"
$ws.Range("C83").Value = "<viza_2>"
$ws.Range("C84").Value = "<viza_3>"
$ws.Range("C89").Value = "<doer>"

# Row 89 previously had a taller custom height to fit the long Russian
# "Исп.:" text; the shorter placeholder goes back to the sheet's default.
$ws.Rows.Item(89).RowHeight = $ws.Rows.Item(88).RowHeight

# --- Restore the scroll position / selection the author left the sheet in ---
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 67 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
$ws.Range("C85").Select()
